$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column D to Text format first so numeric-looking values
# (e.g. "302.80") are stored as literal strings, matching the source
# inline-string cells, instead of being auto-converted to numbers.
$dRange = $ws.Range('D2:D51')
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.210.92'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '2.358.00'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '302.80'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').Value = '95.73'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').Value = '34.13'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '18.70'
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('E13').Value = '  +3.66%  '
$ws.Range('D14').Value = '6.73'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').Value = '2.724.05'
$ws.Range('E15').Value = '  +2.74%  '
$ws.Range('D16').Value = '2.364.98'
$ws.Range('E16').Value = '  +3.16%  '
$ws.Range('D17').Value = '0.797'
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').Value = '43.199.65'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('E20').Value = '  +4.31%  '
$ws.Range('D21').Value = '0.0₃0890'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = '68.21'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('D23').Value = '235.49'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '2.43'
$ws.Range('E26').Value = '  +1.39%  '
$ws.Range('D27').Value = '24.57'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  +0.88%  '
$ws.Range('D30').Value = '31.34'
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('E33').Value = '  +3.63%  '
$ws.Range('D34').Value = '17.27'
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.84'
$ws.Range('E35').Value = '  +5.30%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '4.39'
$ws.Range('E36').Value = '  -1.18%  '
$ws.Range('D37').Value = '2.31'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').Value = '22.51'
$ws.Range('E39').Value = '  +9.89%  '
$ws.Range('E40').Value = '  +1.85%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = '105.79'
$ws.Range('E42').Value = '  -35.72%  '
$ws.Range('D43').Value = '1.945.16'
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('D44').Value = '0.0279'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('E45').Value = '  +4.28%  '
$ws.Range('D46').Value = '9.54'
$ws.Range('E46').Value = '  -8.76%  '
$ws.Range('D47').Value = '2.74'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').Value = '2.587.19'
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('D49').Value = '52.86'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('E50').Value = '  -4.58%  '
$ws.Range('D51').Value = '72.22'
$ws.Range('E51').Value = '  +1.34%  '

# Restore the default (unstyled) cell format on column D, since setting
# NumberFormat above would otherwise leave a stray style index behind.
$ws.Range('G1').Copy()
$dRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0
